$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 977-978, pushing the existing data
# (and everything below it) down by two rows.
$ws.Rows("977:978").Insert()

# Populate the newly inserted row 977 (Coliflor, Primera, fecha 45041).
$ws.Range("A977").Value = 3
$ws.Range("B977").Value = "Femacal de La Calera"
$ws.Range("C977").Value = "Coquimbo"
$ws.Range("D977").Value = 45041
$ws.Range("D977").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E977").Value = 5
$ws.Range("F977").Value = 100112008
$ws.Range("G977").Value = "Coliflor"
$ws.Range("H977").Value = "Sin especificar"
$ws.Range("I977").Value = "Primera"
$ws.Range("J977").Value = 4300
$ws.Range("K977").Value = 1100
$ws.Range("L977").Value = 1200
$ws.Range("M977").Value = 1149
$ws.Range("N977").Value = "$/unidad"
$ws.Range("O977").Value = "Provincia de Quillota"
$ws.Range("P977").Value = 1149
$ws.Range("Q977").Value = 1
$ws.Range("R977").Value = "Hortaliza"

# Populate the newly inserted row 978 (Coliflor, Segunda, fecha 45041).
$ws.Range("A978").Value = 3
$ws.Range("B978").Value = "Femacal de La Calera"
$ws.Range("C978").Value = "Coquimbo"
$ws.Range("D978").Value = 45041
$ws.Range("D978").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E978").Value = 5
$ws.Range("F978").Value = 100112008
$ws.Range("G978").Value = "Coliflor"
$ws.Range("H978").Value = "Sin especificar"
$ws.Range("I978").Value = "Segunda"
$ws.Range("J978").Value = 1300
$ws.Range("K978").Value = 900
$ws.Range("L978").Value = 900
$ws.Range("M978").Value = 900
$ws.Range("N978").Value = "$/unidad"
$ws.Range("O978").Value = "Provincia de Quillota"
$ws.Range("P978").Value = 900
$ws.Range("Q978").Value = 1
$ws.Range("R978").Value = "Hortaliza"
